# Prezenta Algoritmi - sapt. 8 update:
#  - add new student "Catalina Madalina Parca" (sapt. 8 = 1)
#  - bump sapt. 8 mark from 1 to 2 for 4 students who already had an entry
#  - table keeps being sorted alphabetically by name (col B), so the student
#    list below is written out already in that order
#  - restore the previously active selection cell

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
  @("Alessandro Vereș-Pop", 2, 2, 2, $null, 2, $null, $null, $null),
  @("Amanda Hajdu", 1, 1, 1, 1, 1, $null, $null, $null),
  @("Andra Agud", 1, $null, $null, $null, $null, $null, $null, $null),
  @("Attila Bunta", 1, $null, 2, 2, 2, 1, 2, 2),
  @("Cătălina Mădălina Parca", $null, $null, $null, $null, $null, $null, $null, 1),
  @("Claudiu Druța", 2, 1, $null, 1, 2, 1, 2, $null),
  @("Codruț Avram", 1, 1, 2, 1, 1, $null, 2, $null),
  @("Daniela Cionca (Mărie)", 2, 2, 1, 2, 1, 2, 2, 1),
  @("Delia Negrea", 2, 1, 1, 2, 1, 1, 2, $null),
  @("Denisa Cioban", 1, 1, 1, 1, 1, $null, $null, $null),
  @("Levente Nagy", 2, 2, 2, 2, 2, 2, 2, 2),
  @("Luca Șeicaru", 1, $null, 2, 2, 2, 1, 2, 2),
  @("Mark Pop", 2, 2, 2, 2, 2, 1, $null, $null),
  @("Miriam Bacso", 2, 1, $null, $null, $null, $null, $null, $null),
  @("Paul Dobroțchi", 2, 1, 1, 1, 1, $null, $null, $null),
  @("Raul Andrei", 1, 2, 2, 2, 2, 1, $null, $null),
  @("Răzvan Baroi", 2, 2, 2, 2, 2, 1, $null, 1),
  @("Silvia Naghi", 2, 2, 2, 2, 2, 2, $null, 1),
  @("Sorin Fechete", 1, 1, 1, 1, 1, $null, $null, $null),
  @("Victor Lazăr", 1, 2, 2, 2, 2, 2, 2, 2)
)

# Clear the attendance block first so rows that lose a mark (because a
# different student now occupies that row after the alphabetical re-sort)
# don't keep stale values.
$ws.Range("B3:J22").ClearContents()

for ($i = 0; $i -lt $data.Length; $i++) {
  $row = $data[$i]
  $r = 3 + $i
  for ($j = 0; $j -lt $row.Length; $j++) {
    $val = $row[$j]
    if ($null -ne $val) {
      $ws.Cells.Item($r, 2 + $j).Value = $val
    }
  }
}

# Restore the reported active cell / selection
$ws.Range("M17").Select()
